$d = $word.ActiveDocument

# wdAlignParagraphJustify
$wdAlignParagraphJustify = 3

# Paragraphs 7..13 (1-indexed) are the shaded "FFF1A8" e-mail paragraphs.
# For these: every run/paragraph-mark font size goes from 19 half-points (9.5pt)
# to 24 half-points (12pt), and the paragraph gets justified alignment.
for ($i = 7; $i -le 13; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Size = 12
    $p.Alignment = $wdAlignParagraphJustify
}

# Runs wrapped in a Hyperlink are not touched by the paragraph-level font
# assignment above, so resize them explicitly.
foreach ($h in $d.Hyperlinks) {
    $h.Range.Font.Size = 12
}

# Move the "_GoBack" bookmark out of paragraph 8 ("Bonjour Ali,") and into the
# middle of paragraph 9's final run, splitting the word "prochaine" into
# "p" + "rochaine" (". . . pour la p" | "rochaine fois:").
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$p9 = $d.Paragraphs.Item(9)
$p9Text = $p9.Range.Text
$splitOffset = $p9Text.IndexOf("pour la p") + "pour la p".Length
$splitAbs = $p9.Range.Start + $splitOffset
$splitRange = $d.Range($splitAbs, $splitAbs)
$d.Bookmarks.Add("_GoBack", $splitRange)

# Paragraph 14 ("  - Hai") is the final paragraph: its paragraph-mark size
# grows from 24 to 36 half-points (18pt) while its run text only grows from
# 19 to 24 half-points (12pt); it also gains justified alignment.
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Font.Size = 18
$p14.Alignment = $wdAlignParagraphJustify
$p14TextRange = $d.Range($p14.Range.Start, $p14.Range.End - 1)
$p14TextRange.Font.Size = 12

Write-Host "done"
